$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: ricardo.ccorrea1@... -> teste.001@senacsp.edu.br -----------
$ws.Range("A2").Value = "teste.001@senacsp.edu.br"
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:teste.001@senacsp.edu.br")

# A3:A4 picked up a hyperlink together (e.g. via fill/copy from A2) while
# they still showed the old display text, before being overwritten below.
$ws.Hyperlinks.Add($ws.Range("A3:A4"), "mailto:teste.002@senacsp.edu.br", "", "", "teste.001@senacsp.edu.br")

# --- Row 3: henri.cfernandes@... -> teste.002@senacsp.edu.br -----------
$ws.Range("A3").Value = "teste.002@senacsp.edu.br"
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:teste.002@senacsp.edu.br")

# --- Row 4: paulo.cmmartins@... -> teste.003@senacsp.edu.br ------------
$ws.Range("A4").Value = "teste.003@senacsp.edu.br"
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:teste.003@senacsp.edu.br")

# --- "TSC-2025-NOITE-" -> "TSC-2025-TESTE-" in column D ----------------
$ws.Range("D2").Formula = '="TSC-2025-TESTE-" & C2'
$ws.Range("D3").Formula = '="TSC-2025-TESTE-" & C3'
$ws.Range("D4").Formula = '="TSC-2025-TESTE-" & C4'

# --- Final selection left on D4, matching the saved workbook view ------
$ws.Range("D4").Select()
